# Remove the CAAU container row (row 2) from the sequence sheet.
# This reflects that the download/withdrawal info for this row is no
# longer part of the tracked list, shifting the remaining rows (FFAU,
# TGHU, BMOU) up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
